# "added input to data"
#
# Duplicates the existing ValidLogin sheet into a new InvalidLogin sheet
# (right after it), replaces the second data row with the invalid-login
# credentials, and brings the view state (zoom / selection / active tab)
# in line with the authored workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate "ValidLogin" -> new sheet placed immediately after it, then
# rename it. Copy() (rather than Worksheets.Add()) preserves the original
# sheet's formatting/markup so the new sheet starts out identical.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "InvalidLogin"

# Row 2 on the new sheet becomes the "invalid" credentials.
$ws2.Range("A2").Value = "bhanu"
$ws2.Range("B2").Value = "bhanu123"

# ValidLogin: no longer the selected tab; selection becomes the full A1:B2
# data range, zoom stays 145%.
$ws1.Activate()
$ws1.Range("A1:B2").Select()
$excel.ActiveWindow.Zoom = 145

# InvalidLogin: becomes the active/selected tab, zoomed to 160%, with the
# cursor sitting just below the data in B3 (matches a freshly-typed sheet).
$ws2.Activate()
$ws2.Range("B3").Select()
$excel.ActiveWindow.Zoom = 160
